$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.611.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.224.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.76%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.225.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.72%  "

$ws.Range("E10").Value = "  -2.78%  "

$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.390"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.793.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("E14").Value = "  -3.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.817.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.28%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000159"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.65%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.209.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "414.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.81"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.23%  "

$ws.Range("E25").Value = "  -0.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.203"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.494"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.20%  "

$ws.Range("E28").Value = "  -2.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.57%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").Value = "  -3.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.05%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("E34").Value = "  -2.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.25%  "

$ws.Range("E38").Value = "  -2.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.823.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.06%  "

$ws.Range("E40").Value = "  -2.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.727"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.00%  "

$ws.Range("E46").Value = "  -4.59%  "

$ws.Range("E47").Value = "  -5.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "303.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0262"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.57%  "

$ws.Range("E51").Value = "  -0.87%  "
